$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.759.51"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.045.95"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'227.53"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'60.20"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.377"
$ws.Range("E9").Value = "  -1.70%  "
$ws.Range("D10").Value = "'0.0833"
$ws.Range("E10").Value = "  +2.44%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "2.348.70"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "'14.40"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'21.48"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "'5.51"
$ws.Range("E15").Value = "  +6.23%  "
$ws.Range("D16").Value = "'0.765"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("D17").Value = "2.047.61"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "37.753.06"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").Value = "'69.37"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'222.44"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").Value = "'169.05"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "'9.32"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'18.76"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("E32").Value = "  +8.42%  "
$ws.Range("D33").Value = "'4.38"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").Value = "'4.51"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").Value = "'0.0603"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +2.49%  "
$ws.Range("D37").Value = "'2.35"
$ws.Range("E37").Value = "  +4.48%  "
$ws.Range("D38").Value = "'3.47"
$ws.Range("E38").Value = "  +7.32%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").Value = "'18.29"
$ws.Range("E40").Value = "  +9.20%  "
$ws.Range("D41").Value = "1.531.64"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'97.82"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").Value = "'0.0216"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D45").Value = "'4.15"
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("D46").Value = "'0.0890"
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("D49").Value = "'2.94"
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("D50").Value = "'7.02"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "2.237.32"
$ws.Range("E51").Value = "  +0.84%  "
